$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rename the last column from "useBehavioral" (Yes/No flag) to
# "Study Number" (1/2 numeric group id) and rewrite its values.
# ---------------------------------------------------------------------
$ws.Range("L1").Value = "Study Number"

$studyNumberByRow = @{
    2  = 1;  3  = 1;  4  = 1;  5  = 1;  6  = 1;  7  = 1;  8  = 1;  9  = 1
    10 = 1;  11 = 1;  12 = 1;  13 = 1;  14 = 1;  15 = 1;  16 = 1;  17 = 1
    18 = 1;  19 = 1;  20 = 1
    21 = 2;  22 = 2;  23 = 2;  24 = 2;  25 = 2;  26 = 2
}

foreach ($row in $studyNumberByRow.Keys) {
    $ws.Cells.Item($row, 12).Value = $studyNumberByRow[$row]
}

# ---------------------------------------------------------------------
# Fill in previously-blank "Cap Size" values for a few subjects.
# ---------------------------------------------------------------------
$ws.Range("K23").Value = "M"
$ws.Range("K25").Value = "M"
$ws.Range("K26").Value = "M"

# ---------------------------------------------------------------------
# Corrected session counts for a few subjects.
# ---------------------------------------------------------------------
$ws.Range("E23").Value = 3
$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 3

# ---------------------------------------------------------------------
# Update the view: scroll right so column D is the left-most visible
# column, and move the active selection to L23.
# ---------------------------------------------------------------------
[void]$ws.Range("L23").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
